{"js": "// 1) Merge the split run + remove the stray \"_GoBack\" bookmark around\n//    \"Vince will fix curr|ent input read script so reads 1 vs multiple subjects\".\nconst body = context.document.body;\n\nconst fixResults = body.search(\"Vince will fix current input read script so reads 1 vs multiple subjects\", { matchCase: true, ignorePunct: false });\nfixResults.load(\"items\");\nawait context.sync();\n\nif (fixResults.items.length > 0) {\n  // Re-write the paragraph's full range with a single clean run; this merges\n  // the two runs back into one and drops the bookmark that sat between them.\n  const para = fixResults.items[0].paragraphs.getFirst();\n  const rng = para.getRange();\n  rng.insertText(\"Vince will fix current input read script so reads 1 vs multiple subjects\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Add two new bullet points after the \"In paper, will write analysis...\"\n//    bullet, before the blank paragraph that precedes \"John Data Dump of Issues\".\nconst anchorResults = body.search(\"In paper, will write analysis section generically and provide sample scripts in AFNI\", { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n  const anchorPara = anchorResults.items[0].paragraphs.getFirst();\n\n  const p1 = anchorPara.insertParagraph(\"Need to check how short SD pairs are being used\", Word.InsertLocation.after);\n  await context.sync();\n\n  const p2 = p1.insertParagraph(\"SWITCH TO WEIGHTED MEAN FOR GLM\\u2026\", Word.InsertLocation.after);\n  await context.sync();\n\n  // Place the \"_GoBack\" bookmark right at the end of the new final run, matching\n  // where Word leaves it after the user's last edit.\n  const textSearch = p2.search(\"SWITCH TO WEIGHTED MEAN FOR GLM\\u2026\", { matchCase: true });\n  textSearch.load(\"items\");\n  await context.sync();\n\n  if (textSearch.items.length > 0) {\n    const endRange = textSearch.items[0].getRange(\"End\");\n    endRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-Range([string]$text) {\n    $rng = $d.Content\n    $f = $rng.Find\n    $f.ClearFormatting()\n    $f.Forward = $true\n    $f.Wrap = 1\n    $f.MatchCase = $true\n    $f.MatchWholeWord = $false\n    $f.Text = $text\n    if ($f.Execute()) { return $rng }\n    return $null\n}\n\n# ---------------------------------------------------------------------------\n# 1) \"Vince will fix curr|ent input read script ...\" was split across two\n#    runs with a stray \"_GoBack\" bookmark sitting between them. Merge it back\n#    into a single run and drop the bookmark.\n# ---------------------------------------------------------------------------\n$rng1 = Find-Range \"Vince will fix curr\"\nif ($rng1 -ne $null) {\n    $rng1.Expand(4)              # wdParagraph -> whole paragraph incl. the trailing mark\n    $rng1.MoveEnd(1, -1)         # wdCharacter -> exclude the paragraph mark itself\n    $rng1.Delete()\n    $rng1.InsertBefore(\"Vince will fix current input read script so reads 1 vs multiple subjects\")\n}\n\n# ---------------------------------------------------------------------------\n# 2) Add two new bullets after \"In paper, will write analysis section\n#    generically and provide sample scripts in AFNI\", before the blank\n#    paragraph that precedes \"John Data Dump of Issues\".\n# ---------------------------------------------------------------------------\n$rngAnchor = Find-Range \"In paper, will write analysis section generically and provide sample scripts in AFNI\"\nif ($rngAnchor -ne $null) {\n    $rngAnchor.Expand(4)\n    $rngAnchor.InsertParagraphAfter()\n    $rngAnchor.Collapse(0)       # wdCollapseEnd -> start of the freshly-inserted blank paragraph\n    $rngAnchor.Expand(4)\n    $rngAnchor.MoveEnd(1, -1)\n    $rngAnchor.InsertBefore(\"Need to check how short SD pairs are being used\")\n}\n\n$rngBullet2 = Find-Range \"Need to check how short SD pairs are being used\"\nif ($rngBullet2 -ne $null) {\n    $rngBullet2.Expand(4)\n    $rngBullet2.InsertParagraphAfter()\n    $rngBullet2.Collapse(0)\n    $rngBullet2.Expand(4)\n    $rngBullet2.MoveEnd(1, -1)\n\n    $finalText = \"SWITCH TO WEIGHTED MEAN FOR GLM\" + [char]0x2026\n    # Insert with a one-character sentinel tail so the \"_GoBack\" bookmark can be\n    # anchored strictly between the text and the paragraph mark (a bookmark\n    # collapsed exactly at the paragraph's last text offset does not reattach\n    # correctly), then trim the sentinel back off.\n    $rngBullet2.InsertBefore($finalText + \"X\")\n\n    $bmRng = Find-Range ($finalText + \"X\")\n    $bmRng.MoveStart(1, $finalText.Length)\n    $bmRng.Collapse(1)           # wdCollapseStart -> right after the real text, before \"X\"\n    $d.Bookmarks.Add(\"_GoBack\", $bmRng)\n\n    $sentinelRng = Find-Range ($finalText + \"X\")\n    $sentinelRng.MoveStart(1, $finalText.Length)\n    $sentinelRng.Delete()\n}\n"}
